$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 9034
$ws.Range("F3").Value = 1979
$ws.Range("F4").Value = 6622
$ws.Range("F6").Value = 2145
$ws.Range("F7").Value = 603
$ws.Range("F10").Value = 71
$ws.Range("F13").Value = 7
$ws.Range("F14").Value = 84
$ws.Range("F15").Value = 23
$ws.Range("F16").Value = 8943
$ws.Range("F19").Value = 204
$ws.Range("F21").Value = 1845
$ws.Range("F23").Value = 15
$ws.Range("F25").Value = 96
$ws.Range("F28").Value = 1041
$ws.Range("F29").Value = 17
$ws.Range("F31").Value = 559
$ws.Range("F32").Value = 30
$ws.Range("F33").Value = 36
$ws.Range("F34").Value = 546
$ws.Range("F35").Value = 2347
$ws.Range("F36").Value = 879
$ws.Range("F37").Value = 547
$ws.Range("F41").Value = 299
$ws.Range("F43").Value = 9
$ws.Range("F44").Value = 1069
$ws.Range("F46").Value = 23
$ws.Range("F47").Value = 83
$ws.Range("F48").Value = 14
$ws.Range("F49").Value = 4002

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F12").Value = 14
$ws.Range("F14").Value = 17
$ws.Range("F25").Value = 69

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 2343
$ws.Range("F3").Value = 725
$ws.Range("F4").Value = 340
$ws.Range("F5").Value = 14

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2343
$ws.Range("F3").Value = 9035
$ws.Range("F4").Value = 340
$ws.Range("F5").Value = 1979
$ws.Range("F6").Value = 6622
$ws.Range("F8").Value = 2145
$ws.Range("F11").Value = 603
$ws.Range("F15").Value = 71
$ws.Range("F17").Value = 14
$ws.Range("F18").Value = 84
$ws.Range("F19").Value = 8943
$ws.Range("F22").Value = 204
$ws.Range("F23").Value = 1845
$ws.Range("F25").Value = 15
$ws.Range("F26").Value = 96
$ws.Range("F28").Value = 1041
$ws.Range("F29").Value = 17
$ws.Range("F32").Value = 559
$ws.Range("F33").Value = 30
$ws.Range("F34").Value = 36
$ws.Range("F35").Value = 546
$ws.Range("F36").Value = 2347
$ws.Range("F37").Value = 879
$ws.Range("F38").Value = 17
$ws.Range("F40").Value = 547
$ws.Range("F41").Value = 299
$ws.Range("F43").Value = 83
$ws.Range("F44").Value = 4002
$ws.Range("F45").Value = 69
